# This script re-derives the "Tire Type" pipeline numbers after the
# upstream per-signal-column values on Step1_Data changed (the commit's
# "Tire Type Filtering" / extraction cleanup). Step1_Data holds raw
# per-signal fractional values; Step2_Sj is the row-wise running
# cumulative sum of Step1_Data; each Step3_DataPts_<threshold> sheet
# finds, per row, the first column where the Step2_Sj cumulative value
# reaches that row's threshold (col B) and records the 0-based index
# before it (D), the cumulative value there (F) and D-1 (G).

$wb = $excel.ActiveWorkbook

$dataSheetName = "Step1_Data"
$sjSheetName   = "Step2_Sj"
$firstDataCol  = 2   # column B
$lastDataCol   = 36  # column AJ
$dataRows      = @(2, 3, 4, 5, 6)

# --- 1. Apply the new raw values onto Step1_Data -------------------------
$updates = @(
    @{Row=2; Col=5; Val=0.1237204684060884},
    @{Row=2; Col=7; Val=0.4456220838801876},
    @{Row=2; Col=8; Val=0.1554202831739369},
    @{Row=2; Col=10; Val=0.01685728373759012},
    @{Row=2; Col=11; Val=0.02704785073528512},
    @{Row=2; Col=12; Val=0.08083482865667004},
    @{Row=2; Col=13; Val=0.01669810175818519},
    @{Row=2; Col=14; Val=0.07700646557899593},
    @{Row=2; Col=18; Val=0.007509608052988468},
    @{Row=2; Col=20; Val=0.02546205982046139},
    @{Row=2; Col=28; Val=0.02382096619961072},
    @{Row=3; Col=5; Val=0.2787814649009919},
    @{Row=3; Col=7; Val=0.379346472336011},
    @{Row=3; Col=8; Val=0.05424362000115154},
    @{Row=3; Col=10; Val=0.00582739201085851},
    @{Row=3; Col=11; Val=0.07002686958089174},
    @{Row=3; Col=12; Val=0.03801578027060478},
    @{Row=3; Col=13; Val=0.05076961985217018},
    @{Row=3; Col=14; Val=0.03374242251043722},
    @{Row=3; Col=15; Val=0.01687755034502457},
    @{Row=3; Col=18; Val=0.01748273514782861},
    @{Row=3; Col=20; Val=0.005656569561145065},
    @{Row=3; Col=22; Val=0.01183638337962496},
    @{Row=3; Col=24; Val=0.007449212475846418},
    @{Row=3; Col=25; Val=0.01427905087743306},
    @{Row=3; Col=28; Val=0.01566485674998008},
    @{Row=4; Col=4; Val=0.001714494602450862},
    @{Row=4; Col=5; Val=0.2551265672837331},
    @{Row=4; Col=6; Val=0.2154288682764892},
    @{Row=4; Col=7; Val=0.2411969333968549},
    @{Row=4; Col=11; Val=0.09639365453309245},
    @{Row=4; Col=12; Val=0.02040602887623578},
    @{Row=4; Col=13; Val=0.09424572932767201},
    @{Row=4; Col=15; Val=0.01123145274165671},
    @{Row=4; Col=18; Val=0.01577449076726763},
    @{Row=4; Col=21; Val=0.01530881820383216},
    @{Row=4; Col=24; Val=0.02230426136995108},
    @{Row=4; Col=25; Val=0.006720623105210893},
    @{Row=4; Col=27; Val=0.001426484502504124},
    @{Row=4; Col=28; Val=0.002721593013048926},
    @{Row=5; Col=5; Val=0.06703230525802782},
    @{Row=5; Col=6; Val=0.1023731790348799},
    @{Row=5; Col=7; Val=0.2794569729697763},
    @{Row=5; Col=8; Val=0.2728390212949439},
    @{Row=5; Col=12; Val=0.09433318493493846},
    @{Row=5; Col=13; Val=0.04305210067836614},
    @{Row=5; Col=14; Val=0.1122594509362769},
    @{Row=5; Col=17; Val=0.001142684794423614},
    @{Row=5; Col=20; Val=0.004859654624338767},
    @{Row=5; Col=22; Val=0.004083201743710534},
    @{Row=5; Col=25; Val=0.01237868771192732},
    @{Row=5; Col=26; Val=0.00240964062922797},
    @{Row=5; Col=28; Val=0.002539006728323168},
    @{Row=5; Col=29; Val=0.00124090866083899},
    @{Row=6; Col=5; Val=0.337741437284348},
    @{Row=6; Col=6; Val=0.01726272788531873},
    @{Row=6; Col=7; Val=0.3608615291586841},
    @{Row=6; Col=8; Val=0.0272853508974509},
    @{Row=6; Col=10; Val=0.003397632191675581},
    @{Row=6; Col=11; Val=0.02325093387103951},
    @{Row=6; Col=12; Val=0.1180361709232703},
    @{Row=6; Col=13; Val=0.05982802414850919},
    @{Row=6; Col=16; Val=0.001406457561595064},
    @{Row=6; Col=17; Val=0.002042495262617118},
    @{Row=6; Col=18; Val=0.0009723263879604468},
    @{Row=6; Col=19; Val=0.01472920304425142},
    @{Row=6; Col=21; Val=0.005523426988508213},
    @{Row=6; Col=24; Val=0.01591163796458753},
    @{Row=6; Col=28; Val=0.01175064643018398}
)

$wsData = $wb.Worksheets.Item($dataSheetName)
foreach ($u in $updates) {
    $wsData.Cells.Item($u.Row, $u.Col).Value2 = $u.Val
}

# --- 2. Recompute Step2_Sj as the running row-wise cumulative sum --------
$wsSj = $wb.Worksheets.Item($sjSheetName)
foreach ($r in $dataRows) {
    $running = 0.0
    for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
        $v = [double]$wsData.Cells.Item($r, $col).Value2
        $running = $running + $v
        $wsSj.Cells.Item($r, $col).Value2 = $running
    }
}

# --- 3. Recompute each Step3_DataPts_<threshold> sheet --------------------
$thresholdSheets = @(
    "Step3_DataPts_0.5",
    "Step3_DataPts_0.7",
    "Step3_DataPts_0.8",
    "Step3_DataPts_0.9"
)

$colD = 4  # Point_Exceeds_Index
$colF = 6  # Point_Exceeds_Cumulative_Value
$colG = 7  # Pulse_Width

foreach ($sheetName in $thresholdSheets) {
    $wsThresh = $wb.Worksheets.Item($sheetName)
    foreach ($r in $dataRows) {
        $threshold = [double]$wsThresh.Cells.Item($r, 2).Value2

        $foundCol = -1
        $foundVal = 0.0
        for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
            $cum = [double]$wsSj.Cells.Item($r, $col).Value2
            if ($cum -ge $threshold) {
                $foundCol = $col
                $foundVal = $cum
                break
            }
        }

        if ($foundCol -ge 0) {
            $wsThresh.Cells.Item($r, $colD).Value2 = ($foundCol - 1)
            $wsThresh.Cells.Item($r, $colF).Value2 = $foundVal
            $wsThresh.Cells.Item($r, $colG).Value2 = ($foundCol - 2)
        }
    }
}
